# Remove the "Parte III - Feedback Qualitativo (Perguntas Abertas)" block
# (heading, intro sentence, the 4 open questions, and the leading blank
# paragraph that precedes the heading) that sat right after the second
# table in the document, while leaving the blank paragraph and the
# closing "Muito obrigado..." paragraph that follow it untouched.

$d = $word.ActiveDocument

# Locate the heading paragraph and the last question paragraph by text so
# the script is resilient to any paragraph-index quirks; use their Range
# boundaries to build the deletion span.

$startText = "Parte III - Feedback Qualitativo (Perguntas Abertas)"
$endText   = "4. H"

$startRange = $d.Content
$startRange.Find.Execute($startText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingStart = $startRange.Start

$endRange = $d.Content
$endRange.Find.Execute($endText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
# Expand to the full paragraph that contains the match, so the whole
# question (and its trailing paragraph mark) is included in the span.
$endPara = $endRange.Paragraphs(1)
$questionEnd = $endPara.Range.End

# The blank paragraph immediately before the heading must also go: walk
# back over the paragraph mark that precedes the heading paragraph.
$blankRange = $d.Range($headingStart, $headingStart)
$blankPara = $blankRange.Paragraphs(1).Previous()
$deleteStart = $blankPara.Range.Start

$d.Range($deleteStart, $questionEnd).Delete()
